# Update stimulus labels in the trial sheet:
#  - column L ("correct_ans") holds short codes that should be spelled out:
#      b -> center, y -> left, r -> right
#  - any stimulus filename using the old "face" category is renamed to the
#    new "book" category (face//face_NN.jpg -> book//book_NN.jpg), which can
#    show up in columns A-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$ansMap = @{ "b" = "center"; "y" = "left"; "r" = "right" }

for ($r = 2; $r -le $lastRow; $r++) {

    # Columns A-D can contain "face//face_NN.jpg" stimulus paths -> rename to "book"
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $text = $cell.Text
        if ($text -like "*face*") {
            $cell.Value = $text.Replace("face", "book")
        }
    }

    # Column L holds the abbreviated correct-answer code -> spell it out
    $ansCell = $ws.Cells.Item($r, 12)
    $code = $ansCell.Text
    if ($ansMap.ContainsKey($code)) {
        $ansCell.Value = $ansMap[$code]
    }
}
